# Generate Report for Handback
# Update the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps to reflect a fresh report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-07 07:31:42"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H2").Value = "2016-09-07 07:31:36"
$wsZhCn.Range("K2").Value = "2016-09-07 07:31:56"

# --- de-de sheet ---
# H2 ("Correspond Handoff Datetime") shares its timestamp value with
# Overview!G2 ("Latest HO Xliff Generate Date") -- both move to 07:31:42.
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H2").Value = "2016-09-07 07:31:42"
$wsDeDe.Range("K2").Value = "2016-09-07 07:32:16"
